$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# G2 and H2
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# G3 standalone formula
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

# G4:G11 formulas (fill down pattern)
$ws.Range("G4").Formula = "=(D4-D3)*B4/100"
$ws.Range("G5").Formula = "=(D5-D4)*B5/100"
$ws.Range("G6").Formula = "=(D6-D5)*B6/100"
$ws.Range("G7").Formula = "=(D7-D6)*B7/100"
$ws.Range("G8").Formula = "=(D8-D7)*B8/100"
$ws.Range("G9").Formula = "=(D9-D8)*B9/100"
$ws.Range("G10").Formula = "=(D10-D9)*B10/100"
$ws.Range("G11").Formula = "=(D11-D10)*B11/100"

[void]$ws.Range("G1").Select()
